$d = $word.ActiveDocument

# --- 1. Paragraph 1: "...puerto." -> "...puerto y que pueda reservar."
#     Locate the trailing "puerto." and insert the new phrase right before
#     its final period (i.e. right after "puerto").
$findRange = $d.Content.Duplicate
[void]$findRange.Find.Execute("puerto.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$periodPos = $findRange.End - 1
$insertRange = $d.Range($periodPos, $periodPos)
$insertRange.Text = " y que pueda reservar"

# --- 2. Paragraph 2 ("Cual es la informacion mas relevante...") is replaced
#     entirely by a new question. Delete the whole paragraph (including its
#     mark, so the stray proofErr spell-check tags go with it) and retype
#     fresh content further down.
$p2 = $d.Paragraphs(2)
$p2Range = $d.Range($p2.Range.Start, $p2.Range.End)
$p2Range.Delete()

# --- 3. Append the new questions as brand-new paragraphs at the end of the
#     document, plus two trailing blank paragraphs. Using a literal carriage
#     return ("`r") followed by the text on a fresh Range anchored at the
#     current end of story creates a clean new <w:p> for each one.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`rCual es la diferencia entre plaza transito a plaza base"

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`r¿Hay tripulantes tanto en plazas base y transitos?"

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`rCual es la función y definición de un tripulante"

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`rLas funciones de Xunta y Policia/Aduanas solo pueden consultar datos? O también modifican datos"

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`r"

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "`r"
